$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin "Price" values are stored as plain text in this sheet (several
# contain two decimal points, e.g. "39.835.92", or other non-numeric
# formatting quirks) - Excel automatically re-interprets a plain
# assignment such as $range.Value = "86.09" as the *number* 86.09,
# which silently drops meaningful trailing zeros (e.g. "15.90" -> 15.9)
# and changes the stored cell type. Routing the literal text through a
# =T("...") formula and then collapsing it to a static value via
# Copy/PasteSpecial(values) keeps the cell genuinely text-typed without
# touching the cells NumberFormat/style in any way.
function Set-TextValue($range, $text) {
    $escaped = $text.Replace("`"", "`"`"")
    $range.Formula = "=T(`"" + $escaped + "`")"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "39.835.92"
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").Value = "2.218.67"
$ws.Range("E3").Value = "  +1.51%  "

Set-TextValue $ws.Range("D5") "292.69"
$ws.Range("E5").Value = "  -1.29%  "

Set-TextValue $ws.Range("D6") "86.09"
$ws.Range("E6").Value = "  +5.74%  "

Set-TextValue $ws.Range("D7") "0.515"
$ws.Range("E7").Value = "  +1.33%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +1.80%  "

Set-TextValue $ws.Range("D10") "30.76"
$ws.Range("E10").Value = "  +6.19%  "

$ws.Range("E11").Value = "  +2.40%  "

Set-TextValue $ws.Range("D12") "47.34"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("E13").Value = "  +1.71%  "

$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").Value = "2.558.06"
$ws.Range("E15").Value = "  +1.14%  "

Set-TextValue $ws.Range("D16") "14.05"
$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").Value = "2.225.62"
$ws.Range("E17").Value = "  +1.47%  "

$ws.Range("E18").Value = "  +3.55%  "

$ws.Range("D19").Value = "39.802.82"
$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("E20").Value = "  +1.83%  "

Set-TextValue $ws.Range("D21") "11.06"
$ws.Range("E21").Value = "  +8.42%  "

Set-TextValue $ws.Range("D22") "5.81"
$ws.Range("E22").Value = "  +2.30%  "

Set-TextValue $ws.Range("D23") "65.58"
$ws.Range("E23").Value = "  +1.32%  "

Set-TextValue $ws.Range("D24") "236.24"
$ws.Range("E24").Value = "  +4.99%  "

$ws.Range("E25").Value = "  +0.04%  "

Set-TextValue $ws.Range("D26") "2.47"
$ws.Range("E26").Value = "  +2.90%  "

$ws.Range("E27").Value = "  +3.08%  "

Set-TextValue $ws.Range("D28") "22.77"
$ws.Range("E28").Value = "  +1.50%  "

Set-TextValue $ws.Range("D29") "2.11"
$ws.Range("E29").Value = "  -2.19%  "

Set-TextValue $ws.Range("D30") "9.23"
$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("E31").Value = "  +4.51%  "

Set-TextValue $ws.Range("D32") "151.47"
$ws.Range("E32").Value = "  +1.70%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("E34").Value = "  +2.95%  "

Set-TextValue $ws.Range("D35") "0.0719"
$ws.Range("E35").Value = "  +4.35%  "

$ws.Range("E37").Value = "  +7.61%  "

$ws.Range("E38").Value = "  +2.36%  "

Set-TextValue $ws.Range("D39") "15.90"
$ws.Range("E39").Value = "  +5.03%  "

$ws.Range("E40").Value = "  +3.11%  "

$ws.Range("E41").Value = "  +4.94%  "

$ws.Range("E42").Value = "  +5.69%  "

$ws.Range("D43").Value = "2.064.95"
$ws.Range("E43").Value = "  +9.31%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D44") "0.0268"
$ws.Range("E44").Value = "  +4.44%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "17.88"
$ws.Range("E45").Value = "  +12.15%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D46") "2.10"
$ws.Range("E46").Value = "  +1.38%  "

Set-TextValue $ws.Range("D47") "9.95"
$ws.Range("E47").Value = "  +11.62%  "

$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "2.433.97"
$ws.Range("E49").Value = "  +1.36%  "

Set-TextValue $ws.Range("D50") "71.59"
$ws.Range("E50").Value = "  +1.26%  "

Set-TextValue $ws.Range("D51") "88.99"
$ws.Range("E51").Value = "  +3.00%  "

$excel.CutCopyMode = $false

